$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-colour the "PASSED" fill (style s="4") from red (FFC7CE) to green
#    (C6EFCE). Style index 4 already exists (used by the pre-existing
#    I23 "FAILED" cell) so we mutate its Interior colour in place - this
#    updates the shared fill definition rather than creating a brand new
#    unused style for every cell that references it.
# ---------------------------------------------------------------------------
$passInt = $ws.Range("I23").Interior
$passInt.PatternColor = 13561798   # RGB(198,239,206) -> 0xC6EFCE (BGR-encoded long)
$passInt.Color = 13561798

# ---------------------------------------------------------------------------
# 2. All the existing test rows (2-33) had their Actual Result / Test Status
#    parameterised: "Test not executed"/"Not Run" becomes
#    "Signup functionality verified"/"PASSED", and the one row that used to
#    report a failure (row 23) is flipped to a pass as well. The PASSED
#    cells pick up the (now green) style used by column I.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 33; $r++) {
    $ws.Range("H$r").Value = "Signup functionality verified"
    $ws.Range("I$r").Value = "PASSED"

    if ($r -ne 23) {
        $ws.Range("I23").Copy()
        $ws.Range("I$r").PasteSpecial(-4122)
    }
}

# ---------------------------------------------------------------------------
# 3. Insert a brand-new test case row at position 34 ("Verify visibility of
#    resend OTP button on verification page."), pushing the former row 34
#    (TC_033 - resend OTP countdown) down to row 35.
# ---------------------------------------------------------------------------
$ws.Rows(34).Insert()
$ws.Rows(34).RowHeight = 80

# Clone formatting (borders/alignment/fill) for the new row from row 33.
$ws.Range("A33:I33").Copy()
$ws.Range("A34:I34").PasteSpecial(-4122)

$ws.Range("A34").Value = "'33"
$ws.Range("B34").Value = "TC_033"
$ws.Range("C34").Value = "Verify visibility of resend OTP button on verification page."
$ws.Range("D34").Value = "User is on the signup page"
$ws.Range("E34").Value = "Signup test data"
$ws.Range("F34").Value = "1. Navigate to signup page`n2. Perform required actions`n3. Verify expected behavior"
$ws.Range("G34").Value = "Signup functionality should work as expected"
$ws.Range("H34").Value = "Signup functionality verified"
$ws.Range("I34").Value = "PASSED"

# Re-apply the clean (non quote-prefixed) number format/style to A34 and
# give I34 the green PASSED style.
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("I23").Copy()
$ws.Range("I34").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. The row that got pushed down to 35 keeps its original objective text
#    but its serial number / test case id are renumbered (33 -> 34,
#    TC_033 -> TC_034) and its result is parameterised the same way as the
#    rest of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = "'34"
$ws.Range("A33").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("B35").Value = "TC_034"
$ws.Range("H35").Value = "Signup functionality verified"
$ws.Range("I35").Value = "PASSED"
$ws.Range("I23").Copy()
$ws.Range("I35").PasteSpecial(-4122)

$excel.CutCopyMode = $false
